$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "Ergänzungen" ---
$ws.Range("E1").Value = "Ergänzungen"
$ws.Range("E2").Value = "Preis bei 80 Desktops"
$ws.Range("E4").Value = "Preis bei 32 Simultanverbindungen"

# Give column E an explicit custom width (closest achievable to the
# authored 28.265625 produced by real Excel's font-metric autofit).
$ws.Columns.Item(5).ColumnWidth = 27.5

# --- Row 6: PhonerLite compatibility text simplified ---
$ws.Range("B6").Value = "Windows"

# --- New row 8: SessionTalk ---
$ws.Range("A8").Value = "SessionTalk"
$ws.Range("B8").Value = "iOS, Android"
$ws.Range("C8").Value = "Freeware"
$ws.Range("D8").Value = "-"

# --- View state: selection + zoom ---
[void]$ws.Activate()
[void]$ws.Range("C7").Select()
$excel.ActiveWindow.Zoom = 157
